$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("D2").Value = '25.778.45'
$ws.Range("E2").Value = '  -0.31%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("D3").Value = '1.635.66'
$ws.Range("E3").Value = '  -0.10%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("D5").Value = '''215.47'
$ws.Range("E5").Value = '  +0.06%  '

$ws.Range("B6").Value = 'XRP'
$ws.Range("D6").Value = '''0.505'
$ws.Range("E6").Value = '  -0.04%  '

$ws.Range("B7").Value = 'USDC'
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  -0.10%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("D8").Value = '''0.258'
$ws.Range("E8").Value = '  +0.27%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("D9").Value = '''0.0641'
$ws.Range("E9").Value = '  -0.17%  '

$ws.Range("B10").Value = 'Solana'
$ws.Range("D10").Value = '''19.85'
$ws.Range("E10").Value = '  +0.26%  '

$ws.Range("B11").Value = 'TRON'
$ws.Range("D11").Value = '''0.0779'
$ws.Range("E11").Value = '  -0.09%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.668.90'
$ws.Range("E12").Value = '  +1.90%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''4.26'
$ws.Range("E13").Value = '  -0.76%  '

$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("D14").Value = '1.860.96'
$ws.Range("E14").Value = '  -0.15%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("D15").Value = '''0.557'
$ws.Range("E15").Value = '  -0.56%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("D16").Value = '0.0₃0777'
$ws.Range("E16").Value = '  +2.00%  '

$ws.Range("B17").Value = 'Litecoin'
$ws.Range("D17").Value = '''63.22'
$ws.Range("E17").Value = '  +0.38%  '

$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("D18").Value = '25.806.88'
$ws.Range("E18").Value = '  -0.27%  '

$ws.Range("B19").Value = 'Dai'
$ws.Range("D19").Value = '''1.00'
$ws.Range("E19").Value = '  -0.14%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("D20").Value = '''4.44'
$ws.Range("E20").Value = '  +2.60%  '

$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("D21").Value = '''194.33'
$ws.Range("E21").Value = '  -0.15%  '

$ws.Range("B22").Value = 'Avalanche'
$ws.Range("D22").Value = '''9.97'
$ws.Range("E22").Value = '  +0.93%  '

$ws.Range("B23").Value = 'Chainlink'
$ws.Range("D23").Value = '''6.16'
$ws.Range("E23").Value = '  +0.97%  '

$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("D24").Value = '''1.00'
$ws.Range("E24").Value = '  -0.14%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("D25").Value = '''1.76'
$ws.Range("E25").Value = '  -0.70%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("D26").Value = '''139.98'
$ws.Range("E26").Value = '  -0.31%  '

$ws.Range("B27").Value = 'Stellar'
$ws.Range("D27").Value = '''0.121'
$ws.Range("E27").Value = '  -4.09%  '

$ws.Range("B28").Value = 'Cosmos'
$ws.Range("D28").Value = '''6.86'
$ws.Range("E28").Value = '  +0.46%  '

$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("D29").Value = '''15.64'
$ws.Range("E29").Value = '  +1.35%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("D30").Value = '''1.24'
$ws.Range("E30").Value = '  -0.07%  '

$ws.Range("B31").Value = 'Hedera'
$ws.Range("D31").Value = '''0.0493'
$ws.Range("E31").Value = '  +0.68%  '

$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("D32").Value = '''3.35'
$ws.Range("E32").Value = '  +1.70%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("D33").Value = '''3.27'
$ws.Range("E33").Value = '  +1.37%  '

$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("D34").Value = '''1.59'
$ws.Range("E34").Value = '  +1.54%  '

$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("D35").Value = '''2.39'
$ws.Range("E35").Value = '  +0.26%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("D36").Value = '''0.898'
$ws.Range("E36").Value = '  -0.57%  '

$ws.Range("B37").Value = 'MXToken'
$ws.Range("D37").Value = '''2.58'
$ws.Range("E37").Value = '  +0.16%  '

$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("D38").Value = '''0.554'
$ws.Range("E38").Value = '  +0.27%  '

$ws.Range("B39").Value = 'Maker'
$ws.Range("D39").Value = '1.107.73'
$ws.Range("E39").Value = '  -1.68%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("D40").Value = '''0.0157'
$ws.Range("E40").Value = '  +0.52%  '

$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("D41").Value = '''1.00'
$ws.Range("E41").Value = '  +0.40%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("D42").Value = '''5.57'
$ws.Range("E42").Value = '  +0.45%  '

$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("D43").Value = '''0.804'
$ws.Range("E43").Value = '  +0.40%  '

$ws.Range("B44").Value = 'Quant'
$ws.Range("D44").Value = '''99.24'
$ws.Range("E44").Value = '  +1.13%  '

$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("D45").Value = '0.0₆0109'
$ws.Range("E45").Value = '  -4.98%  '

$ws.Range("B46").Value = 'Aave'
$ws.Range("D46").Value = '''55.26'
$ws.Range("E46").Value = '  -0.28%  '

$ws.Range("B47").Value = 'SynthetixNetwork'
$ws.Range("D47").Value = '''2.49'
$ws.Range("E47").Value = '  +12.80%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''7.73'
$ws.Range("E48").Value = '  +0.14%  '

$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").Value = '''0.418'
$ws.Range("E49").Value = '  -2.06%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '''0.0504'
$ws.Range("E50").Value = '  -0.04%  '

$ws.Range("B51").Value = 'Frax'
$ws.Range("D51").Value = '''1.00'
$ws.Range("E51").Value = '  -0.14%  '
